$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 107
$ws.Range("I2").Value = 83.333336
$ws.Range("K2").Value = 83.333336
$ws.Range("M2").Value = 29.666664

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2148.4
$ws.Range("J17").Value = 2148.4
$ws.Range("L17").Value = 6445.200000000001
$ws.Range("N17").Value = -6781.200000000001

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 750
$ws.Range("I18").Value = 750
$ws.Range("K18").Value = 750
$ws.Range("M18").Value = -466

# ALC row 20
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 933.3333
$ws.Range("I20").Value = 933.3333
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 933.3333
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -703.3333
$ws.Range("N20").ClearContents()

# ALC row 35
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 933.3333
$ws.Range("I35").Value = 933.3333
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 933.3333
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -554.3333
$ws.Range("N35").ClearContents()

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1581.9231
$ws.Range("I58").Value = 93.75
$ws.Range("J58").Value = 2243.3333
$ws.Range("K58").Value = 281.25
$ws.Range("L58").Value = 6729.999899999999
$ws.Range("M58").Value = -131.25
$ws.Range("N58").Value = -7029.999899999999

# ALC row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 67
$ws.Range("I99").Value = 67
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 201
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 1297
$ws.Range("N99").ClearContents()

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 14857
$ws.Range("I116").Value = 3995
$ws.Range("J116").Value = 17572.5
$ws.Range("K116").Value = 3995
$ws.Range("L116").Value = 17572.5
$ws.Range("M116").Value = -553
$ws.Range("N116").Value = -24456.5

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 14677.5
$ws.Range("J132").Value = 3800
$ws.Range("L132").Value = 11400
$ws.Range("N132").Value = -16460

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 69.666664
$ws.Range("I5").Value = 60
$ws.Range("K5").Value = 60
$ws.Range("M5").Value = 52

# ARM row 8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

# ARM row 34
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 29999
$ws.Range("I34").Value = 29999
$ws.Range("K34").Value = 29999
$ws.Range("M34").Value = -29728

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 69.666664
$ws.Range("I4").Value = 60
$ws.Range("K4").Value = 60
$ws.Range("M4").Value = 55

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1495
$ws.Range("I86").Value = 1498.3334
$ws.Range("K86").Value = 1498.3334
$ws.Range("M86").Value = -375.3334

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1495
$ws.Range("I89").Value = 1498.3334
$ws.Range("K89").Value = 7491.666999999999
$ws.Range("M89").Value = -1875.666999999999

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 980.1539
$ws.Range("I107").Value = 1182
$ws.Range("J107").Value = 657.2
$ws.Range("K107").Value = 1182
$ws.Range("L107").Value = 657.2
$ws.Range("M107").Value = 738
$ws.Range("N107").Value = -4497.2

# CRP row 2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 8578.714
$ws.Range("I2").Value = 9910.200000000001
$ws.Range("J2").Value = 5250
$ws.Range("K2").Value = 9910.200000000001
$ws.Range("L2").Value = 5250
$ws.Range("M2").Value = -9797.200000000001
$ws.Range("N2").Value = -5476

# CRP row 14
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 3172.7273
$ws.Range("I14").Value = 2766.6667
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 2766.6667
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -2596.6667
$ws.Range("N14").Value = -5340

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3479
$ws.Range("I31").Value = 1747.5
$ws.Range("K31").Value = 1747.5
$ws.Range("M31").Value = -1452.5

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3479
$ws.Range("I34").Value = 1747.5
$ws.Range("K34").Value = 1747.5
$ws.Range("M34").Value = -1545.5

# CRP row 42
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 10056
$ws.Range("I42").Value = 10056
$ws.Range("K42").Value = 10056
$ws.Range("M42").Value = -9463

# CRP row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 17334.875
$ws.Range("I88").Value = 60000
$ws.Range("J88").Value = 11239.857
$ws.Range("K88").Value = 60000
$ws.Range("L88").Value = 11239.857
$ws.Range("M88").Value = -59594
$ws.Range("N88").Value = -12051.857

# CRP row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 17334.875
$ws.Range("I91").Value = 60000
$ws.Range("J91").Value = 11239.857
$ws.Range("K91").Value = 60000
$ws.Range("L91").Value = 11239.857
$ws.Range("M91").Value = -58596
$ws.Range("N91").Value = -14047.857

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3173.625
$ws.Range("I99").Value = 1914.8334
$ws.Range("K99").Value = 1914.8334
$ws.Range("M99").Value = -416.8334

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3173.625
$ws.Range("I126").Value = 1914.8334
$ws.Range("K126").Value = 5744.5002
$ws.Range("M126").Value = -3274.5002

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5066.3335
$ws.Range("I134").Value = 5066.3335
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 15199.0005
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -12664.0005
$ws.Range("N134").ClearContents()

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I7").Value = 49.75
$ws.Range("J7").Value = 125
$ws.Range("K7").Value = 149.25
$ws.Range("L7").Value = 375
$ws.Range("M7").Value = -37.25
$ws.Range("N7").Value = -599

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 400
$ws.Range("I11").Value = 100
$ws.Range("K11").Value = 300
$ws.Range("M11").Value = -160

# CUL row 25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 491.41666
$ws.Range("I25").Value = 187.75
$ws.Range("J25").Value = 1098.75
$ws.Range("K25").Value = 563.25
$ws.Range("L25").Value = 3296.25
$ws.Range("M25").Value = -394.25
$ws.Range("N25").Value = -3634.25

# CUL row 30
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 491.41666
$ws.Range("I30").Value = 187.75
$ws.Range("J30").Value = 1098.75
$ws.Range("K30").Value = 563.25
$ws.Range("L30").Value = 3296.25
$ws.Range("M30").Value = -461.25
$ws.Range("N30").Value = -3500.25

# CUL row 43
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 1000
$ws.Range("J43").Value = 1000
$ws.Range("L43").Value = 3000
$ws.Range("N43").Value = -3228

# CUL row 44
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 784.7143
$ws.Range("I44").Value = 201
$ws.Range("J44").Value = 1222.5
$ws.Range("K44").Value = 603
$ws.Range("L44").Value = 3667.5
$ws.Range("M44").Value = -205
$ws.Range("N44").Value = -4463.5

# CUL row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 783
$ws.Range("I99").Value = 783
$ws.Range("K99").Value = 2349
$ws.Range("M99").Value = -103

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2282.353
$ws.Range("J131").Value = 2282.353
$ws.Range("L131").Value = 6847.059
$ws.Range("N131").Value = -16927.059

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6502097
$ws.Range("I11").Value = 13666667
$ws.Range("J11").Value = 3431567.2
$ws.Range("K11").Value = 13666667
$ws.Range("L11").Value = 3431567.2
$ws.Range("M11").Value = -13666528
$ws.Range("N11").Value = -3431845.2

# GSM row 12
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 1950
$ws.Range("I12").Value = 1950
$ws.Range("K12").Value = 1950
$ws.Range("M12").Value = -1810

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6352.846
$ws.Range("J70").Value = 6486.125
$ws.Range("L70").Value = 6486.125
$ws.Range("N70").Value = -7026.125

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6352.846
$ws.Range("J73").Value = 6486.125
$ws.Range("L73").Value = 6486.125
$ws.Range("N73").Value = -8358.125

# GSM row 116
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3308.6
$ws.Range("I126").Value = 3510.75
$ws.Range("K126").Value = 10532.25
$ws.Range("M126").Value = -8062.25

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 10000
$ws.Range("I16").Value = 10000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -9830
$ws.Range("N16").ClearContents()

# LTW row 92
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# WVR row 20
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 2002499.5
$ws.Range("I20").Value = 4999
$ws.Range("K20").Value = 4999
$ws.Range("M20").Value = -4759

# WVR row 127
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 45000
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
